$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("465:465").Insert()

$ws.Range("A465").Value = 8
$ws.Range("B465").Value = "Terminal La Palmera de La Serena"
$ws.Range("C465").Value = "Coquimbo"
$ws.Range("D465").Value = 44585
$ws.Range("E465").Value = 4
$ws.Range("F465").Value = 100112024
$ws.Range("G465").Value = "Choclo"
$ws.Range("H465").Value = "Dulce o Americano"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 40000
$ws.Range("K465").Value = 150
$ws.Range("L465").Value = 200
$ws.Range("M465").Value = 175
$ws.Range("N465").Value = "$/unidad"
$ws.Range("O465").Value = "Provincia del Elquí"
$ws.Range("P465").Value = 175
$ws.Range("Q465").Value = 1
$ws.Range("R465").Value = "Hortaliza"
